# Workbook was re-uploaded with a new "class3/image3" column pair (E:F)
# appended, mostly mirroring the existing C:D "class2/image2" column pair,
# plus a couple of touch-ups to C:D for rows 14-21 (new class boundary
# shifted, and the last group of rows lost its C:D values altogether).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header pair: E1 = "class3", F1 = "image3" ---------------------
$ws.Range("E1").Value = "class3"
$ws.Range("F1").Value = "image3"

# --- 2. Touch-ups to the existing C:D columns (rows 14-21) -----------------
$ws.Range("C14").Value = "car"
$ws.Range("D14").Value = "图像3"
$ws.Range("D15").Value = "图像3"
$ws.Range("D16").Value = "图像3"
$ws.Range("D17").Value = "图像3"

$ws.Range("C18:D21").ClearContents()

# --- 3. Populate the new E:F columns for rows 2-21 --------------------------
$dataEF = @(
    @("car","图像2"),
    @("car","图像2"),
    @("car","图像2"),
    @("car","图像2"),
    @("car","图像2"),
    @("motorcycle","图像2"),
    @("motorcycle","图像2"),
    @("truck","图像2"),
    @("car","图像1"),
    @("car","图像1"),
    @("car","图像1"),
    @("car","图像1"),
    @("bus","图像4"),
    @("car","图像4"),
    @("car","图像4"),
    @("car","图像4"),
    @("car","图像3"),
    @("car","图像3"),
    @("car","图像3"),
    @("car","图像3")
)

for ($i = 0; $i -lt $dataEF.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $dataEF[$i][0]
    $ws.Cells.Item($row, 6).Value = $dataEF[$i][1]
}

# --- 4. Match the new selection recorded in the saved file ------------------
$ws.Range("F1").Select() | Out-Null
